# Scheduled-runner style refresh of market/profit columns (H:N) across several
# sheets. Mirrors the upstream diff: some rows lose their stale price data
# (cleared entirely), some get fresh price data inserted where none existed,
# and some get their existing price data refreshed to new values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALC: rows 125-141 lost all their currentAveragePrice/Leve* columns (H:N).
# ---------------------------------------------------------------------------
$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("H125:N141").ClearContents()

# ---------------------------------------------------------------------------
# CRP: rows 62, 65, 132, 134 got refreshed price data.
# ---------------------------------------------------------------------------
$wsCRP = $wb.Worksheets.Item("CRP")

$wsCRP.Range("H62").Value = 12300
$wsCRP.Range("I62").Value = 3400
$wsCRP.Range("J62").Value = 16750
$wsCRP.Range("K62").Value = 3400
$wsCRP.Range("L62").Value = 16750
$wsCRP.Range("M62").Value = -2776
$wsCRP.Range("N62").Value = -17998

$wsCRP.Range("H65").Value = 12300
$wsCRP.Range("I65").Value = 3400
$wsCRP.Range("J65").Value = 16750
$wsCRP.Range("K65").Value = 17000
$wsCRP.Range("L65").Value = 83750
$wsCRP.Range("M65").Value = -13880
$wsCRP.Range("N65").Value = -89990

$wsCRP.Range("H132").Value = 336999.66
$wsCRP.Range("I132").Value = 502500
$wsCRP.Range("J132").Value = 5999
$wsCRP.Range("K132").Value = 1507500
$wsCRP.Range("L132").Value = 17997
$wsCRP.Range("M132").Value = -1504970
$wsCRP.Range("N132").Value = -23057

$wsCRP.Range("H134").Value = 8462
$wsCRP.Range("I134").Value = 9578
$wsCRP.Range("J134").Value = 3998
$wsCRP.Range("K134").Value = 28734
$wsCRP.Range("L134").Value = 11994
$wsCRP.Range("M134").Value = -26199
$wsCRP.Range("N134").Value = -17064

# ---------------------------------------------------------------------------
# CUL: row 16 reset to zeros (and dropped LeveProfitHQ / N16); rows 120-141
# (minus 135, already populated) gained fresh price data.
# ---------------------------------------------------------------------------
$wsCUL = $wb.Worksheets.Item("CUL")

$wsCUL.Range("H16").Value = 0
$wsCUL.Range("I16").Value = 0
$wsCUL.Range("J16").Value = 0
$wsCUL.Range("K16").Value = 0
$wsCUL.Range("L16").Value = 0
$wsCUL.Range("N16").ClearContents()

$wsCUL.Range("H120").Value = 0
$wsCUL.Range("I120").Value = 0
$wsCUL.Range("J120").Value = 0
$wsCUL.Range("K120").Value = 0
$wsCUL.Range("L120").Value = 0

$wsCUL.Range("H121").Value = 750
$wsCUL.Range("I121").Value = 0
$wsCUL.Range("J121").Value = 750
$wsCUL.Range("K121").Value = 0
$wsCUL.Range("L121").Value = 2250
$wsCUL.Range("N121").Value = -4870

$wsCUL.Range("H122").Value = 1938.5
$wsCUL.Range("I122").Value = 1591.8
$wsCUL.Range("J122").Value = 2516.3333
$wsCUL.Range("K122").Value = 14326.2
$wsCUL.Range("L122").Value = 22646.9997
$wsCUL.Range("M122").Value = -11876.2
$wsCUL.Range("N122").Value = -27546.9997

$wsCUL.Range("H123").Value = 5000
$wsCUL.Range("I123").Value = 0
$wsCUL.Range("J123").Value = 5000
$wsCUL.Range("K123").Value = 0
$wsCUL.Range("L123").Value = 15000
$wsCUL.Range("N123").Value = -19900

$wsCUL.Range("H124").Value = 4747.4546
$wsCUL.Range("I124").Value = 0
$wsCUL.Range("J124").Value = 4747.4546
$wsCUL.Range("K124").Value = 0
$wsCUL.Range("L124").Value = 14242.3638
$wsCUL.Range("N124").Value = -24062.3638

$wsCUL.Range("H125").Value = 0
$wsCUL.Range("I125").Value = 0
$wsCUL.Range("J125").Value = 0
$wsCUL.Range("K125").Value = 0
$wsCUL.Range("L125").Value = 0

$wsCUL.Range("H126").Value = 4550
$wsCUL.Range("I126").Value = 2780
$wsCUL.Range("J126").Value = 7500
$wsCUL.Range("K126").Value = 8340
$wsCUL.Range("L126").Value = 22500
$wsCUL.Range("M126").Value = -3400
$wsCUL.Range("N126").Value = -32380

$wsCUL.Range("H127").Value = 0
$wsCUL.Range("I127").Value = 0
$wsCUL.Range("J127").Value = 0
$wsCUL.Range("K127").Value = 0
$wsCUL.Range("L127").Value = 0

$wsCUL.Range("H128").Value = 0
$wsCUL.Range("I128").Value = 0
$wsCUL.Range("J128").Value = 0
$wsCUL.Range("K128").Value = 0
$wsCUL.Range("L128").Value = 0

$wsCUL.Range("H129").Value = 1300
$wsCUL.Range("I129").Value = 1000
$wsCUL.Range("J129").Value = 2500
$wsCUL.Range("K129").Value = 3000
$wsCUL.Range("L129").Value = 7500
$wsCUL.Range("M129").Value = 2000
$wsCUL.Range("N129").Value = -17500

$wsCUL.Range("H130").Value = 1000
$wsCUL.Range("I130").Value = 1000
$wsCUL.Range("J130").Value = 0
$wsCUL.Range("K130").Value = 3000
$wsCUL.Range("L130").Value = 0
$wsCUL.Range("M130").Value = 2020

$wsCUL.Range("H131").Value = 1515.2307
$wsCUL.Range("I131").Value = 1266.6666
$wsCUL.Range("J131").Value = 1589.8
$wsCUL.Range("K131").Value = 3799.9998
$wsCUL.Range("L131").Value = 4769.4
$wsCUL.Range("M131").Value = 1240.0002
$wsCUL.Range("N131").Value = -14849.4

$wsCUL.Range("H132").Value = 45000
$wsCUL.Range("I132").Value = 0
$wsCUL.Range("J132").Value = 45000
$wsCUL.Range("K132").Value = 0
$wsCUL.Range("L132").Value = 405000
$wsCUL.Range("N132").Value = -410060

$wsCUL.Range("H133").Value = 50
$wsCUL.Range("I133").Value = 50
$wsCUL.Range("J133").Value = 0
$wsCUL.Range("K133").Value = 150
$wsCUL.Range("L133").Value = 0
$wsCUL.Range("M133").Value = 4910

$wsCUL.Range("H134").Value = 2000
$wsCUL.Range("I134").Value = 2000
$wsCUL.Range("J134").Value = 0
$wsCUL.Range("K134").Value = 6000
$wsCUL.Range("L134").Value = 0
$wsCUL.Range("M134").Value = -930

$wsCUL.Range("H136").Value = 0
$wsCUL.Range("I136").Value = 0
$wsCUL.Range("J136").Value = 0
$wsCUL.Range("K136").Value = 0
$wsCUL.Range("L136").Value = 0

$wsCUL.Range("H137").Value = 7499.75
$wsCUL.Range("I137").Value = 15000
$wsCUL.Range("J137").Value = 4999.6665
$wsCUL.Range("K137").Value = 45000
$wsCUL.Range("L137").Value = 14998.9995
$wsCUL.Range("M137").Value = -39900
$wsCUL.Range("N137").Value = -25198.9995

$wsCUL.Range("H138").Value = 3032
$wsCUL.Range("I138").Value = 0
$wsCUL.Range("J138").Value = 3032
$wsCUL.Range("K138").Value = 0
$wsCUL.Range("L138").Value = 9096
$wsCUL.Range("N138").Value = -19376

$wsCUL.Range("H139").Value = 1000000
$wsCUL.Range("I139").Value = 1000000
$wsCUL.Range("J139").Value = 0
$wsCUL.Range("K139").Value = 3000000
$wsCUL.Range("L139").Value = 0
$wsCUL.Range("M139").Value = -2994860

$wsCUL.Range("H140").Value = 9000
$wsCUL.Range("I140").Value = 9000
$wsCUL.Range("J140").Value = 0
$wsCUL.Range("K140").Value = 27000
$wsCUL.Range("L140").Value = 0
$wsCUL.Range("M140").Value = -21820

$wsCUL.Range("H141").Value = 50750
$wsCUL.Range("I141").Value = 100000
$wsCUL.Range("J141").Value = 1500
$wsCUL.Range("K141").Value = 300000
$wsCUL.Range("L141").Value = 4500
$wsCUL.Range("M141").Value = -294820
$wsCUL.Range("N141").Value = -14860

# ---------------------------------------------------------------------------
# GSM: row 130 refreshed; LeveProfitNQ (M130) dropped.
# ---------------------------------------------------------------------------
$wsGSM = $wb.Worksheets.Item("GSM")
$wsGSM.Range("H130").Value = 69000
$wsGSM.Range("I130").Value = 0
$wsGSM.Range("J130").Value = 69000
$wsGSM.Range("K130").Value = 0
$wsGSM.Range("L130").Value = 69000
$wsGSM.Range("M130").ClearContents()
$wsGSM.Range("N130").Value = -79040

# ---------------------------------------------------------------------------
# WVR: row 136 refreshed.
# ---------------------------------------------------------------------------
$wsWVR = $wb.Worksheets.Item("WVR")
$wsWVR.Range("H136").Value = 2381.3076
$wsWVR.Range("I136").Value = 2195.7
$wsWVR.Range("K136").Value = 6587.099999999999
$wsWVR.Range("M136").Value = -4037.099999999999
